# Apply the chart-related edits described by the commit: resize/reposition
# the embedded line chart's graphic frame, and make the value axis use an
# explicit (unlinked) "General" number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# --- Resize the chart's anchor ("to" cell moves from M24-ish to Q37-ish) ---
# The anchor "from" corner is unchanged; only the bottom-right ("to") corner
# moves, which is equivalent to growing the chart's Width/Height while
# keeping Left/Top fixed.
$co.Width = 720.75
$co.Height = 506

# --- Value axis: add an explicit, unlinked "General" number format ---
$valAx = $chart.Axes(2)
$valAx.NumberFormatLinked = 0
$valAx.NumberFormat = "General"
